$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: the "type" metadata row — Provincia (C) and Mes y año (D) move from
# measure to dimension, matching refArea's classification.
$ws.Range("C3").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "iaest-dimension:mes-y-ano"

# Row 4: the "kind" row — Provincia (C) and Mes y año (D) columns now flagged
# as dimensions ("dim") instead of measures ("medida").
$ws.Range("C4").Value = "dim"
$ws.Range("D4").Value = "dim"

# Row 5: the "uri/format" row — Provincia column gets its own URI type.
$ws.Range("C5").Value = "URI-Provincia"
